$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# --- New row 196 ---
$ws.Range("B196").Value = 6
$ws.Range("C196").Value = "Autenticação e autorização com tokens JWT"
$ws.Range("D196").Value = 76
$ws.Range("E196").Value = "Restrição de conteúdo: cliente só recupera seus pedidos"
$ws.Range("F196").Value = " `n0:07`n6. Autenticação e autorização com tokens JWT`n76. Restrição de conteúdo: cliente só recupera seus pedidos`naula com ATUALIZAÇÃO"
$ws.Range("G196").Value = "`n`n`n"

# --- New row 197 ---
$ws.Range("B197").Value = 6
$ws.Range("C197").Value = "Autenticação e autorização com tokens JWT"
$ws.Range("D197").Value = 76
$ws.Range("E197").Value = "Restrição de conteúdo: cliente só recupera seus pedidos"
$ws.Range("F197").Value = "5:10`n6. Autenticação e autorização com tokens JWT`n76. Restrição de conteúdo: cliente só recupera seus pedidos`ncriação de endpoint para buscar os pedidos do cliente que esta logado"
$ws.Range("G197").Value = ""

# --- New row 198 ---
$ws.Range("B198").Value = 6
$ws.Range("C198").Value = "Autenticação e autorização com tokens JWT"
$ws.Range("D198").Value = 76
$ws.Range("E198").Value = "Restrição de conteúdo: cliente só recupera seus pedidos"
$ws.Range("F198").Value = "5:51`n6. Autenticação e autorização com tokens JWT`n76. Restrição de conteúdo: cliente só recupera seus pedidos`nna criação do endpoint findPage não foi inserido o atributo value da anotação @RequestMapping.. pois foi reaproveitado o endpoint principal ""/pedidos"""
$ws.Range("G198").Value = ""

# Copy formatting from row above (row 195) into the new rows
$ws.Range("B195:G195").Copy()
$ws.Range("B196:B198").PasteSpecial(-4122)
$ws.Range("C196:C198").PasteSpecial(-4122)
$ws.Range("D196:D198").PasteSpecial(-4122)
$ws.Range("E196:E198").PasteSpecial(-4122)
$ws.Range("F196:F198").PasteSpecial(-4122)
$ws.Range("G196:G198").PasteSpecial(-4122)

$ws.Rows.Item(196).RowHeight = 75
$ws.Rows.Item(197).RowHeight = 75
$ws.Rows.Item(198).RowHeight = 90

# Update sheet view to match target
$ws.Application.ActiveWindow.ScrollRow = 193
$ws.Range("C197").Select()

# Update workbook window size
$excel.ActiveWindow.Width = 23730
$excel.ActiveWindow.Height = 7785
